# Changes of 1st april 2022
# Updates ShipmentTrackNum (col C) and, where present, PackageTrackNum (col D)
# for rows 2-22 on Sheet1 with a new batch of tracking numbers, mirroring the
# shared-strings additions in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> new tracking number text.
$trackNums = @{
    2  = "320018208097"
    3  = "320018208101"
    4  = "320018208134"
    5  = "320018208156"
    6  = "320018208190"
    7  = "320018208215"
    8  = "320018208248"
    9  = "320018208260"
    10 = "320018208292"
    11 = "320018208318"
    12 = "320018208351"
    13 = "320018208373"
    14 = "320018208400"
    15 = "320018208421"
    16 = "320018208454"
    17 = "320018208476"
    18 = "320018208513"
    19 = "320018208535"
    20 = "320018208568"
    21 = "320018208580"
    22 = "320018208616"
}

# Rows where column D mirrors column C (PackageTrackNum present).
$dRows = @(5, 6, 7, 13, 14, 15, 16, 17)

for ($row = 2; $row -le 22; $row++) {
    $newVal = $trackNums[$row]

    $cols = @("C")
    if ($dRows -contains $row) {
        $cols += "D"
    }

    foreach ($col in $cols) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)

        # Write the digits-only text as a formula returning a text constant so
        # Excel treats it as text (not a number), then convert the formula to
        # a plain cached value via copy / paste-values. This keeps the cell's
        # type as a shared string without touching its number format/style.
        $cell.Formula = '="' + $newVal + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = 0
